$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = $true
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = $false
$ws.Range("I2").Value = $false
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = $true
$ws.Range("L2").Value = $false
$ws.Range("M2").Value = $true
$ws.Range("N2").Value = $false
$ws.Range("O2").Value = $false
$ws.Range("P2").Value = $false

# Row 3
$ws.Range("B3").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = $true
$ws.Range("H3").Value = $false
$ws.Range("I3").Value = $false
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = $true

# Row 4
$ws.Range("E4").Value = $true
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = $false
$ws.Range("I4").Value = $false
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = $true
$ws.Range("L4").Value = $false
$ws.Range("M4").Value = $true
$ws.Range("N4").Value = $false
$ws.Range("O4").Value = $false
$ws.Range("P4").Value = $false

# Row 5
$ws.Range("C5").Value = $true
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = $true
$ws.Range("H5").Value = $false
$ws.Range("I5").Value = $true
$ws.Range("J5").Value = $false
$ws.Range("L5").Value = $false
$ws.Range("M5").Value = $false
$ws.Range("O5").Value = $false
$ws.Range("P5").Value = $false

# Row 6
$ws.Range("B6").Value = $false
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = $true
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = $false
$ws.Range("I6").Value = $false
$ws.Range("J6").Value = $false
$ws.Range("K6").Value = $true
$ws.Range("L6").Value = $false
$ws.Range("M6").Value = $true
$ws.Range("O6").Value = $false
$ws.Range("P6").Value = $false

# Row 7
$ws.Range("E7").Value = $true
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = $false
$ws.Range("I7").Value = $false
$ws.Range("J7").Value = $false
$ws.Range("K7").Value = $true
$ws.Range("L7").Value = $false
$ws.Range("M7").Value = $true
$ws.Range("N7").Value = $false
$ws.Range("O7").Value = $false
$ws.Range("P7").Value = $false

# Row 8
$ws.Range("C8").Value = $true
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("H8").Value = $false
$ws.Range("I8").Value = $true
$ws.Range("J8").Value = $false
$ws.Range("L8").Value = $false
$ws.Range("M8").Value = $false
$ws.Range("O8").Value = $false
$ws.Range("P8").Value = $false

# Row 9
$ws.Range("J9").Value = $false

# Row 10
$ws.Range("E10").Value = $true
$ws.Range("J10").Value = $false

# Row 11
$ws.Range("B11").Value = $true
$ws.Range("H11").Value = $false
$ws.Range("J11").Value = $false
$ws.Range("K11").Value = $false

# Row 12
$ws.Range("E12").Value = $true
$ws.Range("G12").Value = $true
$ws.Range("H12").Value = $false
$ws.Range("I12").Value = $false
$ws.Range("J12").Value = $false
$ws.Range("K12").Value = $true
$ws.Range("L12").Value = $false
$ws.Range("M12").Value = $true
$ws.Range("N12").Value = $false
$ws.Range("O12").Value = $false
$ws.Range("P12").Value = $false

# Row 13
$ws.Range("B13").Value = $false
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = $false
$ws.Range("J13").Value = $false
$ws.Range("K13").Value = $true
$ws.Range("L13").Value = $true
$ws.Range("N13").Value = $true
$ws.Range("O13").Value = $false
$ws.Range("P13").Value = $false

# Row 14
$ws.Range("B14").Value = $false
$ws.Range("D14").Value = $true
$ws.Range("F14").Value = $true
$ws.Range("G14").Value = $true
$ws.Range("H14").Value = $false
$ws.Range("I14").Value = $false
$ws.Range("J14").Value = $false
$ws.Range("K14").Value = $true
$ws.Range("L14").Value = $false
$ws.Range("M14").Value = $true
$ws.Range("P14").Value = $false

# Row 15
$ws.Range("B15").Value = $false
$ws.Range("H15").Value = $false
$ws.Range("I15").Value = $false
$ws.Range("J15").Value = $false
$ws.Range("K15").Value = $true
$ws.Range("M15").Value = $true
$ws.Range("N15").Value = $true
$ws.Range("P15").Value = $false

# Row 16
$ws.Range("E16").Value = $true
$ws.Range("G16").Value = $true
$ws.Range("H16").Value = $false
$ws.Range("I16").Value = $false
$ws.Range("J16").Value = $false
$ws.Range("K16").Value = $true
$ws.Range("L16").Value = $false
$ws.Range("M16").Value = $true
$ws.Range("N16").Value = $false
$ws.Range("O16").Value = $false
$ws.Range("P16").Value = $false

# Row 17
$ws.Range("B17").Value = $false
$ws.Range("H17").Value = $false
$ws.Range("I17").Value = $false
$ws.Range("J17").Value = $false
$ws.Range("K17").Value = $true
$ws.Range("L17").Value = $true
$ws.Range("N17").Value = $true
$ws.Range("O17").Value = $false
$ws.Range("P17").Value = $false
